# Generate Report for Handback
# Refresh the handoff/handback timestamps recorded in the handback-status report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date for the 0f3b7744 handback file
$overview.Range("G2").Value = "2016-08-30 00:49:58"

# zh-cn sheet: refresh handoff/handback datetimes for the 0f3b7744 file row
$zhcn.Range("H2").Value = "2016-08-30 00:49:53"
$zhcn.Range("K2").Value = "2016-08-30 00:50:18"

# de-de sheet: refresh handoff/handback datetimes for the 0f3b7744 file row
$dede.Range("H2").Value = "2016-08-30 00:49:58"
$dede.Range("K2").Value = "2016-08-30 00:50:25"
